$wb = $excel.ActiveWorkbook

# Opportunity to Engagement validation message list: "Equity Advisors" -> "Creditor Advisors"
$ws1 = $wb.Worksheets.Item("AddOpportunity")
$ws1.Range("C2").Value = "Creditor Advisors"

# Expense/SIC requirement validation text: insert the "Total Debt HL represents(MM)" requirement
$ws3 = $wb.Worksheets.Item("AddContact")
$newSicValidation = @"
Opportunity Detail - SIC Code., Opportunity Description - Opportunity Description., Marketing Information - Client Description., Marketing Information - Total Debt HL represents(MM), input zero if it’s N/A, Marketing Information - Select "Total Debt(MM) Confirmed" checkbox to confirm Total Debt is most up to date., Estimated Fees - Retainer, input zero if there's no Retainer fee., Estimated Fees - Progress/Monthly Fee., Estimated Fees - Contingent Fee., Referral Information - Referral Contact name is required., Referral Information - "Legal Advisor to Company" is required. Please update this field with the correct value, Referral Information - "Legal Advisor to Creditor" is required. Please update this field with the correct value, HL Internal Team - Team must include the following roles: Initiator, Seller, Principal, Manager, Associate(Optional), Analyst(Optional)., Legal Matters - Confidentiality Agreement, Conflicts Check - A Conflicts Check was completed more than 30 days ago. A new Conflicts Check must be completed., Administration - Estimated Closed Date., Administration - "Women Led" is required. Please update this field with the correct value, Administration - Date Engaged - Date of Executed Retainer or similar document., Opportunity Contacts - Add at least one Primary Opportunity Contact., Opportunity Contacts - Add at least one Billing Contact.
"@
$ws3.Range("I2").Value = $newSicValidation
